$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2021.4
$ws.Range("I40").Value = 2230.5
$ws.Range("J40").Value = 1828.3846
$ws.Range("K40").Value = 2230.5
$ws.Range("L40").Value = 1828.3846
$ws.Range("M40").Value = -2055.5
$ws.Range("N40").Value = -2178.3846
$ws.Range("H55").Value = 369.18182
$ws.Range("I55").Value = 307.625
$ws.Range("J55").Value = 533.3333
$ws.Range("K55").Value = 307.625
$ws.Range("L55").Value = 533.3333
$ws.Range("M55").Value = -93.625
$ws.Range("N55").Value = -961.3333
$ws.Range("H98").Value = 5972.55
$ws.Range("I98").Value = 3841.7222
$ws.Range("K98").Value = 3841.7222
$ws.Range("M98").Value = -2343.7222
$ws.Range("H103").Value = 1100
$ws.Range("I103").Value = 300
$ws.Range("J103").Value = 1500
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 4500
$ws.Range("M103").Value = -314
$ws.Range("N103").Value = -5672
$ws.Range("H122").Value = 5972.55
$ws.Range("I122").Value = 3841.7222
$ws.Range("K122").Value = 11525.1666
$ws.Range("M122").Value = -9075.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2620.1567
$ws.Range("I32").Value = 2492.0137
$ws.Range("K32").Value = 2492.0137
$ws.Range("M32").Value = -2205.0137
$ws.Range("H74").Value = 865.5
$ws.Range("I74").Value = 680.61536
$ws.Range("J74").Value = 1666.6666
$ws.Range("K74").Value = 680.61536
$ws.Range("L74").Value = 1666.6666
$ws.Range("M74").Value = 193.38464
$ws.Range("N74").Value = -3414.6666
$ws.Range("H77").Value = 865.5
$ws.Range("I77").Value = 680.61536
$ws.Range("J77").Value = 1666.6666
$ws.Range("K77").Value = 3403.0768
$ws.Range("L77").Value = 8333.333000000001
$ws.Range("M77").Value = 964.9232000000002
$ws.Range("N77").Value = -17069.333
$ws.Range("H97").Value = 800
$ws.Range("I97").Value = 800
$ws.Range("K97").Value = 800
$ws.Range("M97").Value = -304
$ws.Range("H122").Value = 1598.1428
$ws.Range("I122").Value = 1377.4
$ws.Range("J122").Value = 2150
$ws.Range("K122").Value = 4132.200000000001
$ws.Range("L122").Value = 6450
$ws.Range("M122").Value = -1682.200000000001
$ws.Range("N122").Value = -11350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 43277.445
$ws.Range("J132").Value = 43277.445
$ws.Range("L132").Value = 43277.445
$ws.Range("N132").Value = -53397.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1318.9166
$ws.Range("I134").Value = 1188.8182
$ws.Range("K134").Value = 3566.4546
$ws.Range("M134").Value = -1031.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 82.25
$ws.Range("I2").Value = 40
$ws.Range("K2").Value = 240
$ws.Range("M2").Value = -127
$ws.Range("H4").Value = 981935.25
$ws.Range("I4").Value = 143.22223
$ws.Range("J4").Value = 5399999.5
$ws.Range("K4").Value = 429.66669
$ws.Range("L4").Value = 16199998.5
$ws.Range("M4").Value = -317.66669
$ws.Range("N4").Value = -16200222.5
$ws.Range("H7").Value = 466.66666
$ws.Range("I7").Value = 485.42856
$ws.Range("J7").Value = 401
$ws.Range("K7").Value = 1456.28568
$ws.Range("L7").Value = 1203
$ws.Range("M7").Value = -1344.28568
$ws.Range("N7").Value = -1427
$ws.Range("H23").Value = 387.55
$ws.Range("I23").Value = 705
$ws.Range("J23").Value = 251.5
$ws.Range("K23").Value = 2115
$ws.Range("L23").Value = 754.5
$ws.Range("M23").Value = -1880
$ws.Range("N23").Value = -1224.5
$ws.Range("H34").Value = 1888.8889
$ws.Range("I34").Value = 863.25
$ws.Range("J34").Value = 2709.4
$ws.Range("K34").Value = 2589.75
$ws.Range("L34").Value = 8128.200000000001
$ws.Range("M34").Value = -2505.75
$ws.Range("N34").Value = -8296.200000000001
$ws.Range("H39").Value = 2069.739
$ws.Range("J39").Value = 1810.2
$ws.Range("L39").Value = 5430.6
$ws.Range("N39").Value = -6018.6
$ws.Range("H55").Value = 2832.8333
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 3199.4
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 9598.200000000001
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -9952.200000000001
$ws.Range("H131").Value = 13335590
$ws.Range("I131").Value = 76923410
$ws.Range("J131").Value = 2660.742
$ws.Range("K131").Value = 230770230
$ws.Range("L131").Value = 7982.226000000001
$ws.Range("M131").Value = -230765190
$ws.Range("N131").Value = -18062.226
$ws.Range("H139").Value = 1669.6923
$ws.Range("I139").Value = 1737.0952
$ws.Range("J139").Value = 1386.6
$ws.Range("K139").Value = 5211.2856
$ws.Range("L139").Value = 4159.799999999999
$ws.Range("M139").Value = -71.28560000000016
$ws.Range("N139").Value = -14439.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 737.125
$ws.Range("I97").Value = 737.125
$ws.Range("K97").Value = 737.125
$ws.Range("M97").Value = -241.125
$ws.Range("H107").Value = 1053.8572
$ws.Range("I107").Value = 1068.5454
$ws.Range("K107").Value = 1068.5454
$ws.Range("M107").Value = 851.4546
$ws.Range("H122").Value = 1547.4166
$ws.Range("I122").Value = 1705.5555
$ws.Range("K122").Value = 5116.666499999999
$ws.Range("M122").Value = -2666.666499999999
$ws.Range("H132").Value = 1685.9697
$ws.Range("I132").Value = 1486.8422
$ws.Range("J132").Value = 1956.2142
$ws.Range("K132").Value = 4460.5266
$ws.Range("L132").Value = 5868.642599999999
$ws.Range("M132").Value = -1930.5266
$ws.Range("N132").Value = -10928.6426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 730.3333
$ws.Range("I22").Value = 508.16666
$ws.Range("K22").Value = 508.16666
$ws.Range("M22").Value = -213.16666
$ws.Range("H27").Value = 730.3333
$ws.Range("I27").Value = 508.16666
$ws.Range("K27").Value = 508.16666
$ws.Range("M27").Value = -401.16666
$ws.Range("H93").Value = 1234.0769
$ws.Range("I93").Value = 671.44446
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 671.44446
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = 576.55554
$ws.Range("N93").Value = -4996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11305739
$ws.Range("I122").Value = 14445739
$ws.Range("J122").Value = 1736.8
$ws.Range("K122").Value = 43337217
$ws.Range("L122").Value = 5210.4
$ws.Range("M122").Value = -43334767
$ws.Range("N122").Value = -10110.4
